# Update the title cell on the "Inputs and Outputs" sheet from "Inputs" to
# "Results Summary and Inputs" (commit: update proforma in/out sheet title).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inputs and Outputs")

$ws.Range("A1").Value = "Results Summary and Inputs"

# Leave the cursor on A1 (top-left / default cell) for this sheet.
$ws.Range("A1").Select()
